$d = $word.ActiveDocument

# --- 1. Add the new "Abstract Title" paragraph style ---------------------
# styleId ends up "AbstractTitle" (spaces stripped) while NameLocal keeps
# the space, matching w:name val="Abstract Title".
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = "Normal"
$abstractTitle.NextParagraphStyle = "Abstract"
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060   # wdColor BGR for RGB(0x34,0x5A,0x8A) -> w:color 345A8A

# --- 2. Abstract style: tighten space-before from 15pt to 5pt ------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5

# --- 3. ImportTok character style gains green bold colouring -------------
$importTok = $d.Styles.Item("ImportTok")
$importTok.Font.Color = 32768   # wdColor BGR for RGB(0x00,0x80,0x00) -> w:color 008000
$importTok.Font.Bold = $true

# --- 4. BuiltInTok character style gains green colouring -----------------
$builtInTok = $d.Styles.Item("BuiltInTok")
$builtInTok.Font.Color = 32768   # wdColor BGR for RGB(0x00,0x80,0x00) -> w:color 008000
